# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to columns H-N across several rows
# on sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR as per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 6888.5
$ws.Range("I34").Value = 6888.5
$ws.Range("K34").Value = 6888.5
$ws.Range("M34").Value = -6685.5

$ws.Range("H36").Value = 6888.5
$ws.Range("I36").Value = 6888.5
$ws.Range("K36").Value = 6888.5
$ws.Range("M36").Value = -6173.5

$ws.Range("H87").Value = 27250
$ws.Range("J87").Value = 27166.666
$ws.Range("L87").Value = 27166.666
$ws.Range("N87").Value = -29662.666

$ws.Range("H90").Value = 27250
$ws.Range("J90").Value = 27166.666
$ws.Range("L90").Value = 81499.99800000001
$ws.Range("N90").Value = -93979.99800000001

$ws.Range("H112").Value = 1859.45
$ws.Range("I112").Value = 1294.5
$ws.Range("J112").Value = 1922.2222
$ws.Range("K112").Value = 3883.5
$ws.Range("L112").Value = 5766.6666
$ws.Range("M112").Value = -2775.5
$ws.Range("N112").Value = -7982.6666

$ws.Range("H116").Value = 3934.3635
$ws.Range("I116").Value = 3911
$ws.Range("K116").Value = 3911
$ws.Range("M116").Value = -469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 2900
$ws.Range("I19").Value = 2900
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 2900
$ws.Range("L19").Value = $null
$ws.Range("M19").Value = -2671
$ws.Range("N19").Value = 0

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = $null
$ws.Range("N113").Value = 0

$ws.Range("H119").Value = 49666.332
$ws.Range("J119").Value = 49666.332
$ws.Range("L119").Value = 49666.332
$ws.Range("N119").Value = -59342.332

$ws.Range("H122").Value = 3007.88
$ws.Range("I122").Value = 2794.647
$ws.Range("K122").Value = 8383.940999999999
$ws.Range("M122").Value = -5933.940999999999

$ws.Range("H132").Value = 7502.75
$ws.Range("I132").Value = 4998.5
$ws.Range("K132").Value = 14995.5
$ws.Range("M132").Value = -12465.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2452.75
$ws.Range("J86").Value = 1999
$ws.Range("L86").Value = 1999
$ws.Range("N86").Value = -4245

$ws.Range("H89").Value = 2452.75
$ws.Range("J89").Value = 1999
$ws.Range("L89").Value = 9995
$ws.Range("N89").Value = -21227

$ws.Range("H105").Value = 4028.5293
$ws.Range("I105").Value = 2582.5
$ws.Range("J105").Value = 7499
$ws.Range("K105").Value = 2582.5
$ws.Range("L105").Value = 7499
$ws.Range("M105").Value = -835.5
$ws.Range("N105").Value = -10993

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1182.6428
$ws.Range("I16").Value = 1035.1538
$ws.Range("K16").Value = 1035.1538
$ws.Range("M16").Value = -748.1538

$ws.Range("H19").Value = 1667185
$ws.Range("I19").Value = 1667185
$ws.Range("K19").Value = 1667185
$ws.Range("M19").Value = -1667015

$ws.Range("H24").Value = 1667185
$ws.Range("I24").Value = 1667185
$ws.Range("K24").Value = 1667185
$ws.Range("M24").Value = -1667015

$ws.Range("H31").Value = 6110.794
$ws.Range("I31").Value = 3586.9546
$ws.Range("K31").Value = 3586.9546
$ws.Range("M31").Value = -3291.9546

$ws.Range("H34").Value = 6110.794
$ws.Range("I34").Value = 3586.9546
$ws.Range("K34").Value = 3586.9546
$ws.Range("M34").Value = -3384.9546

$ws.Range("H113").Value = 1182.6428
$ws.Range("I113").Value = 1035.1538
$ws.Range("K113").Value = 1035.1538
$ws.Range("M113").Value = 1134.8462

$ws.Range("H115").Value = 70500
$ws.Range("J115").Value = 70500
$ws.Range("L115").Value = 70500
$ws.Range("N115").Value = -72850

$ws.Range("H122").Value = 4140.2856
$ws.Range("I122").Value = 4140.2856
$ws.Range("K122").Value = 12420.8568
$ws.Range("M122").Value = -9970.856800000001

$ws.Range("H132").Value = 972
$ws.Range("I132").Value = 972
$ws.Range("K132").Value = 2916
$ws.Range("M132").Value = -386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 42832.555
$ws.Range("J131").Value = 1850.4762
$ws.Range("L131").Value = 5551.4286
$ws.Range("N131").Value = -15631.4286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 30086.5
$ws.Range("I34").Value = 30000
$ws.Range("J34").Value = 30173
$ws.Range("K34").Value = 30000
$ws.Range("L34").Value = 30173
$ws.Range("M34").Value = -29732
$ws.Range("N34").Value = -30709

$ws.Range("H76").Value = 30086.5
$ws.Range("I76").Value = 30000
$ws.Range("J76").Value = 30173
$ws.Range("K76").Value = 30000
$ws.Range("L76").Value = 30173
$ws.Range("M76").Value = -29685
$ws.Range("N76").Value = -30803

$ws.Range("H79").Value = 30086.5
$ws.Range("I79").Value = 30000
$ws.Range("J79").Value = 30173
$ws.Range("K79").Value = 30000
$ws.Range("L79").Value = 30173
$ws.Range("M79").Value = -28908
$ws.Range("N79").Value = -32357

$ws.Range("H80").Value = 1726799
$ws.Range("J80").Value = 2514395.2
$ws.Range("L80").Value = 2514395.2
$ws.Range("N80").Value = -2516391.2

$ws.Range("H83").Value = 1726799
$ws.Range("J83").Value = 2514395.2
$ws.Range("L83").Value = 12571976
$ws.Range("N83").Value = -12581960

$ws.Range("H107").Value = 1194
$ws.Range("I107").Value = 975.8182
$ws.Range("J107").Value = 1794
$ws.Range("K107").Value = 975.8182
$ws.Range("L107").Value = 1794
$ws.Range("M107").Value = 944.1818
$ws.Range("N107").Value = -5634

$ws.Range("H113").Value = 2250
$ws.Range("I113").Value = 2250
$ws.Range("K113").Value = 2250
$ws.Range("M113").Value = -80

$ws.Range("H122").Value = 2999.1428
$ws.Range("I122").Value = 2999
$ws.Range("K122").Value = 8997
$ws.Range("M122").Value = -6547

$ws.Range("H132").Value = 1073.6
$ws.Range("I132").Value = 854.2308
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 2562.6924
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -32.69239999999991
$ws.Range("N132").Value = -12558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4572.4375
$ws.Range("J40").Value = 4878.636
$ws.Range("L40").Value = 4878.636
$ws.Range("N40").Value = -5150.636

$ws.Range("H61").Value = 1046.8334
$ws.Range("I61").Value = 1046.8334
$ws.Range("K61").Value = 1046.8334
$ws.Range("M61").Value = -844.8334

$ws.Range("H109").Value = 50284.668
$ws.Range("J109").Value = 50284.668
$ws.Range("L109").Value = 50284.668
$ws.Range("N109").Value = -53058.668

$ws.Range("H113").Value = 1046.8334
$ws.Range("I113").Value = 1046.8334
$ws.Range("K113").Value = 1046.8334
$ws.Range("M113").Value = 1123.1666

$ws.Range("H122").Value = 5499.4546
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050

$ws.Range("H129").Value = 60000
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null

$ws.Range("H131").Value = 66999.5
$ws.Range("J131").Value = 66999.5
$ws.Range("L131").Value = 66999.5
$ws.Range("N131").Value = -77079.5

$ws.Range("H132").Value = 2946.1538
$ws.Range("I132").Value = 2516.5
$ws.Range("J132").Value = 3137.111
$ws.Range("K132").Value = 7549.5
$ws.Range("L132").Value = 9411.332999999999
$ws.Range("M132").Value = -5019.5
$ws.Range("N132").Value = -14471.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 39.5
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = $null

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = $null
$ws.Range("N24").Value = 0

$ws.Range("H35").Value = 39.5
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = $null

$ws.Range("H39").Value = 19999.5
$ws.Range("J39").Value = 19999.5
$ws.Range("L39").Value = 19999.5
$ws.Range("N39").Value = -20825.5

$ws.Range("H42").Value = 44999
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = $null

$ws.Range("H119").Value = 79997.75
$ws.Range("J119").Value = 79997.75
$ws.Range("L119").Value = 79997.75
$ws.Range("N119").Value = -89673.75

$ws.Range("H129").Value = 34444.445
$ws.Range("I129").Value = 15000
$ws.Range("K129").Value = 15000
$ws.Range("M129").Value = -10000

$ws.Range("H132").Value = 2115.7307
$ws.Range("I132").Value = 2160.32
$ws.Range("K132").Value = 6480.960000000001
$ws.Range("M132").Value = -3950.960000000001

$ws.Range("H136").Value = 1690.4615
$ws.Range("I136").Value = 1332.7778
$ws.Range("K136").Value = 3998.3334
$ws.Range("M136").Value = -1448.3334
